$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = $ws.Range("B2").Value()
$ws.Range("B2").Value = "Totals"

$ws.Range("B3").Select()
